# iApp.xlsx test-data update ("new changes for oauth, bot detection, file upload"):
#   - new OAuth-style admin credentials (email/password) replacing the old ones,
#     each now hyperlinked, plus a hyperlink on the existing goto URL
#   - new signin JSON payload matching the new credentials (request row shrinks)
#   - new automation steps for the Proposals / Finalizations / Appendix / file
#     upload flow (rows 10-18)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Stash a pristine copy of B2's "Hyperlink" cell format (font only, no extra
# alignment) in a scratch cell. Hyperlinks.Add() always reformats the target
# range, so any cell we hyperlink needs this reapplied afterwards to land on
# the same cellXf the workbook already uses for B2.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# goto URL (B2) becomes a real hyperlink; Excel splits the "#" into the
# address/location pair.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "https://qa.iapp.cool/", "/")

# ---------------------------------------------------------------------------
# Row 4 / Row 5: new login credentials (email + password), each hyperlinked
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = "borhadeashish27+aadmin11@yopmail.com"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:borhadeashish27+aadmin11@yopmail.com")

$ws.Range("Z1").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "Admin@123"
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Admin@123")

# Restore B2/C4/C5 to the original "Hyperlink" cellXf (Hyperlinks.Add just
# bumped all three to a freshly-minted lookalike style).
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# ---------------------------------------------------------------------------
# Row 7: sign-in request now posts the new credentials; row shrinks to fit
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = '{"email": "borhadeashish27+aadmin11@yopmail.com","password": "Admin@123"}'
$ws.Rows.Item(7).RowHeight = 57.6

# ---------------------------------------------------------------------------
# New rows 10-18: Proposals -> Finalizations -> Appendix -> file upload flow
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "waitfortext"
$ws.Range("B10").Value = "Proposals"

$ws.Range("A11").Value = "click"
$ws.Range("B11").Value = "Proposals in sidebar"
$ws.Range("D11").Value = 1000
$ws.Range("E11").Value = 2000

$ws.Range("A12").Value = "click"
$ws.Range("B12").Value = "Finalizations tab "
$ws.Range("D12").Value = 1000
$ws.Range("E12").Value = 2000

$ws.Range("A13").Value = "scrollto"
$ws.Range("B13").Value = "Finalize Proposal"
$ws.Range("D13").Value = 1000
$ws.Range("E13").Value = 2000

$ws.Range("A14").Value = "click"
$ws.Range("B14").Value = "Finalize Proposal"
$ws.Range("D14").Value = 1000
$ws.Range("E14").Value = 2000

# Row 15 previously held only a stray empty formatted cell (B15); reset it to
# the plain default style before filling in the new step.
$ws.Range("A1").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("A15").Value = "click"
$ws.Range("B15").Value = "Appendix in sidebar"
$ws.Range("D15").Value = 1000
$ws.Range("E15").Value = 2000

$ws.Range("A16").Value = "waitfortext"
$ws.Range("B16").Value = "Proposed"
$ws.Range("D16").Value = 1000
$ws.Range("E16").Value = 2000

$ws.Range("A17").Value = "fileupload"
$ws.Range("B17").Value = 'input[accept="application/pdf"]'
$ws.Range("C17").Value = "./uploads/sample.pdf"

$ws.Range("A18").Value = "waitfortext"
$ws.Range("B18").Value = "sample.pdf"

# ---------------------------------------------------------------------------
# Reset the saved cursor position away from the stale A8 selection
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
